$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# Give column B (Unvani) its own best-fit-ish width, independent from column C,
# which keeps its original 22.125 width (splits the old merged B:C col entry).
$ws.Columns("B").ColumnWidth = 31

# --- Fill in "Sicil" numbers down column C (rows 3-99), sequential starting at 1 ---
for ($row = 3; $row -le 99; $row++) {
    $ws.Cells.Item($row, 3).Value = $row - 2
}

# Two cells in column C (C9, C92) had an "odd" style out compared to the rest of the
# column (C9 used the red/alert style, C92 used the header-ish style carried over from
# column A/B on that row). Normalize their formatting to match the rest of column C
# (same style as C3) without touching their neighbours.
$ws.Range("C3").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C92").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection / active cell ---
$ws.Range("E2").Select()
